# Generate Report for Handoff
# Adds two new localization files (653759b1-... and 9309a87c-...) that are
# "Ready for handoff" / "In Translation" to the Overview / zh-cn / de-de
# sheets, pushing the existing ".localization-config" row down.

$wb = $excel.ActiveWorkbook

$newMd1 = "653759b1-0659-4f88-806f-bfbe00d47595.md"
$newMd2 = "9309a87c-4022-404b-85f2-0629e1d1186b.md"

$zhXlf1 = "653759b1-0659-4f88-806f-bfbe00d47595.ff152c00aad098610c7c9e69ff346d5d4610ee3a.zh-cn.xlf"
$zhXlf2 = "9309a87c-4022-404b-85f2-0629e1d1186b.95bc8be20954ffa36113c680f7a4538969322d24.zh-cn.xlf"
$deXlf1 = "653759b1-0659-4f88-806f-bfbe00d47595.ff152c00aad098610c7c9e69ff346d5d4610ee3a.de-de.xlf"
$deXlf2 = "9309a87c-4022-404b-85f2-0629e1d1186b.95bc8be20954ffa36113c680f7a4538969322d24.de-de.xlf"

$zhDate = "2016-02-23 07:09:39"
$deDate = "2016-02-23 07:09:52"

$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/a68f0a7f751b0fef6ac9b515dd556b38333dfbf9/e2e/"
$cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/a68f0a7f751b0fef6ac9b515dd556b38333dfbf9/.localization-config"
$zhHtBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/db65d9bc454d3871961a000d5e76bff34dd2d03a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/"
$deHtBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/746dac158f290ce3e8579847b1b040872c6b7a37/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/"

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Remember what row 4 (".localization-config") looked like so we can move it
# down to make room for the two new rows.
$cfgDisplay = $wsOverview.Range("A4").Value2

$wsOverview.Range("A6").Value = $cfgDisplay
$wsOverview.Range("B6").Value = "Not to be localized"
$wsOverview.Range("C6").Value = "Not to be localized"

$wsOverview.Range("A4").Value = $newMd1
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"

$wsOverview.Range("A5").Value = $newMd2
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), ($mdBase + "082a55fa-75ee-4dea-aed1-abe836dff4ca.md"), [Type]::Missing, [Type]::Missing, "082a55fa-75ee-4dea-aed1-abe836dff4ca.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), ($mdBase + "bbd06bdf-2e3d-4523-a1c2-48a8c127cc6f.md"), [Type]::Missing, [Type]::Missing, "bbd06bdf-2e3d-4523-a1c2-48a8c127cc6f.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), ($mdBase + $newMd1), [Type]::Missing, [Type]::Missing, $newMd1) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), ($mdBase + $newMd2), [Type]::Missing, [Type]::Missing, $newMd2) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), $cfgUrl, [Type]::Missing, [Type]::Missing, $cfgDisplay) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A6").Value = $cfgDisplay
$wsZh.Range("B6").Value = "Not to be localized"
$wsZh.Range("D6").Value = "0001-01-01 00:00:00"
$wsZh.Range("G6").Value = "0001-01-01 00:00:00"
$wsZh.Range("H6").Value = "Ignored"

$wsZh.Range("A4").Value = $newMd1
$wsZh.Range("B4").Value = "Ready for handoff"
$wsZh.Range("C4").Value = $zhXlf1
$wsZh.Range("D4").Value = $zhDate
$wsZh.Range("G4").Value = "0001-01-01 00:00:00"
$wsZh.Range("H4").Value = "Include"

$wsZh.Range("A5").Value = $newMd2
$wsZh.Range("B5").Value = "Ready for handoff"
$wsZh.Range("C5").Value = $zhXlf2
$wsZh.Range("D5").Value = $zhDate
$wsZh.Range("G5").Value = "0001-01-01 00:00:00"
$wsZh.Range("H5").Value = "Include"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), ($mdBase + "082a55fa-75ee-4dea-aed1-abe836dff4ca.md"), [Type]::Missing, [Type]::Missing, "082a55fa-75ee-4dea-aed1-abe836dff4ca.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), ($zhHtBase + "082a55fa-75ee-4dea-aed1-abe836dff4ca.fb38fb412606b3b937c05a6387bb0cdc49a978ba.zh-cn.xlf"), [Type]::Missing, [Type]::Missing, "082a55fa-75ee-4dea-aed1-abe836dff4ca.fb38fb412606b3b937c05a6387bb0cdc49a978ba.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), ($mdBase + "bbd06bdf-2e3d-4523-a1c2-48a8c127cc6f.md"), [Type]::Missing, [Type]::Missing, "bbd06bdf-2e3d-4523-a1c2-48a8c127cc6f.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), ($zhHtBase + "bbd06bdf-2e3d-4523-a1c2-48a8c127cc6f.5eabff815f1236a33edf87df9a041e18ee28c9f5.zh-cn.xlf"), [Type]::Missing, [Type]::Missing, "bbd06bdf-2e3d-4523-a1c2-48a8c127cc6f.5eabff815f1236a33edf87df9a041e18ee28c9f5.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), ($mdBase + $newMd1), [Type]::Missing, [Type]::Missing, $newMd1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C4"), ($zhHtBase + $zhXlf1), [Type]::Missing, [Type]::Missing, $zhXlf1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), ($mdBase + $newMd2), [Type]::Missing, [Type]::Missing, $newMd2) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C5"), ($zhHtBase + $zhXlf2), [Type]::Missing, [Type]::Missing, $zhXlf2) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A6"), $cfgUrl, [Type]::Missing, [Type]::Missing, $cfgDisplay) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A6").Value = $cfgDisplay
$wsDe.Range("B6").Value = "Not to be localized"
$wsDe.Range("D6").Value = "0001-01-01 00:00:00"
$wsDe.Range("G6").Value = "0001-01-01 00:00:00"
$wsDe.Range("H6").Value = "Ignored"

$wsDe.Range("A4").Value = $newMd1
$wsDe.Range("B4").Value = "Ready for handoff"
$wsDe.Range("C4").Value = $deXlf1
$wsDe.Range("D4").Value = $deDate
$wsDe.Range("G4").Value = "0001-01-01 00:00:00"
$wsDe.Range("H4").Value = "Include"

$wsDe.Range("A5").Value = $newMd2
$wsDe.Range("B5").Value = "Ready for handoff"
$wsDe.Range("C5").Value = $deXlf2
$wsDe.Range("D5").Value = $deDate
$wsDe.Range("G5").Value = "0001-01-01 00:00:00"
$wsDe.Range("H5").Value = "Include"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), ($mdBase + "082a55fa-75ee-4dea-aed1-abe836dff4ca.md"), [Type]::Missing, [Type]::Missing, "082a55fa-75ee-4dea-aed1-abe836dff4ca.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), ($deHtBase + "082a55fa-75ee-4dea-aed1-abe836dff4ca.fb38fb412606b3b937c05a6387bb0cdc49a978ba.de-de.xlf"), [Type]::Missing, [Type]::Missing, "082a55fa-75ee-4dea-aed1-abe836dff4ca.fb38fb412606b3b937c05a6387bb0cdc49a978ba.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), ($mdBase + "bbd06bdf-2e3d-4523-a1c2-48a8c127cc6f.md"), [Type]::Missing, [Type]::Missing, "bbd06bdf-2e3d-4523-a1c2-48a8c127cc6f.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), ($deHtBase + "bbd06bdf-2e3d-4523-a1c2-48a8c127cc6f.5eabff815f1236a33edf87df9a041e18ee28c9f5.de-de.xlf"), [Type]::Missing, [Type]::Missing, "bbd06bdf-2e3d-4523-a1c2-48a8c127cc6f.5eabff815f1236a33edf87df9a041e18ee28c9f5.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), ($mdBase + $newMd1), [Type]::Missing, [Type]::Missing, $newMd1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C4"), ($deHtBase + $deXlf1), [Type]::Missing, [Type]::Missing, $deXlf1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), ($mdBase + $newMd2), [Type]::Missing, [Type]::Missing, $newMd2) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C5"), ($deHtBase + $deXlf2), [Type]::Missing, [Type]::Missing, $deXlf2) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A6"), $cfgUrl, [Type]::Missing, [Type]::Missing, $cfgDisplay) | Out-Null
